$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-153) holds the "Förändrad" date serial; update every
# occurrence of 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C153").Value = 45175
